$wb = $excel.ActiveWorkbook

# Move the selection on the previously active sheet (PricingStatus) before
# adding/activating the new sheet, so its saved cursor position matches.
$pricingStatus = $wb.Worksheets.Item("PricingStatus")
$pricingStatus.Select()
$pricingStatus.Range("L34").Select()

# Add the new "MemberStatus" worksheet at the end of the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "MemberStatus"

# Header row
$newSheet.Range("A1").Value = "Id"
$newSheet.Range("B1").Value = "Name"

# Data rows
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "Waiting for the task"

$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "Active"

$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = "On hold"

$newSheet.Range("A5").Value = 4
$newSheet.Range("B5").Value = "Leaver"

# Make the new sheet the active/selected one.
$newSheet.Select()
$newSheet.Range("J17").Select()
